$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1018.75
$ws.Range("I125").Value = 250
$ws.Range("J125").Value = 1787.5
$ws.Range("K125").Value = 2250
$ws.Range("L125").Value = 16087.5
$ws.Range("M125").Value = 210
$ws.Range("N125").Value = -21007.5
$ws.Range("H132").Value = 962386.25
$ws.Range("I132").Value = 1634
$ws.Range("J132").Value = 49000000
$ws.Range("K132").Value = 4902
$ws.Range("L132").Value = 147000000
$ws.Range("M132").Value = -2372
$ws.Range("N132").Value = -147005060
$ws.Range("H137").Value = 2792747
$ws.Range("I137").Value = 7835310
$ws.Range("J137").Value = 1328.125
$ws.Range("K137").Value = 23505930
$ws.Range("L137").Value = 3984.375
$ws.Range("M137").Value = -23503380
$ws.Range("N137").Value = -9084.375
$ws.Range("H138").Value = 2418742.8
$ws.Range("I138").Value = 1147.2
$ws.Range("J138").Value = 3792376.5
$ws.Range("K138").Value = 3441.6
$ws.Range("L138").Value = 11377129.5
$ws.Range("M138").Value = 1698.4
$ws.Range("N138").Value = -11387409.5
$ws.Range("H141").Value = 2618.389
$ws.Range("I141").Value = 2166.5
$ws.Range("J141").Value = 4200
$ws.Range("K141").Value = 6499.5
$ws.Range("L141").Value = 12600
$ws.Range("M141").Value = -1319.5
$ws.Range("N141").Value = -22960

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 45546744
$ws.Range("I61").Value = 58883548
$ws.Range("K61").Value = 58883548
$ws.Range("M61").Value = -58883336
$ws.Range("H132").Value = 61055.383
$ws.Range("I132").Value = 51983.15
$ws.Range("K132").Value = 155949.45
$ws.Range("M132").Value = -153419.45
$ws.Range("H136").Value = 45546744
$ws.Range("I136").Value = 58883548
$ws.Range("K136").Value = 176650644
$ws.Range("M136").Value = -176648094

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 18302.928
$ws.Range("I86").Value = 29750.75
$ws.Range("J86").Value = 3039.1667
$ws.Range("K86").Value = 29750.75
$ws.Range("L86").Value = 3039.1667
$ws.Range("M86").Value = -28627.75
$ws.Range("N86").Value = -5285.1667
$ws.Range("H89").Value = 18302.928
$ws.Range("I89").Value = 29750.75
$ws.Range("J89").Value = 3039.1667
$ws.Range("K89").Value = 148753.75
$ws.Range("L89").Value = 15195.8335
$ws.Range("M89").Value = -143137.75
$ws.Range("N89").Value = -26427.8335
$ws.Range("H134").Value = 5752.8335
$ws.Range("I134").Value = 4498.6816
$ws.Range("J134").Value = 9201.75
$ws.Range("K134").Value = 13496.0448
$ws.Range("L134").Value = 27605.25
$ws.Range("M134").Value = -10961.0448
$ws.Range("N134").Value = -32675.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 27570300
$ws.Range("J58").Value = 1499.6666
$ws.Range("L58").Value = 1499.6666
$ws.Range("N58").Value = -1905.6666
$ws.Range("H127").Value = 31241.666
$ws.Range("J127").Value = 31241.666
$ws.Range("L127").Value = 31241.666
$ws.Range("N127").Value = -41161.666
$ws.Range("H132").Value = 53027.5
$ws.Range("I132").Value = 3637.3333
$ws.Range("J132").Value = 74194.71000000001
$ws.Range("K132").Value = 10911.9999
$ws.Range("L132").Value = 222584.13
$ws.Range("M132").Value = -8381.999899999999
$ws.Range("N132").Value = -227644.13
$ws.Range("H133").Value = 38537.727
$ws.Range("J133").Value = 38537.727
$ws.Range("L133").Value = 38537.727
$ws.Range("N133").Value = -43597.727
$ws.Range("H134").Value = 48361.305
$ws.Range("I134").Value = 1147.625
$ws.Range("J134").Value = 156278.28
$ws.Range("K134").Value = 3442.875
$ws.Range("L134").Value = 468834.84
$ws.Range("M134").Value = -907.875
$ws.Range("N134").Value = -473904.84
$ws.Range("H135").Value = 55632.668
$ws.Range("J135").Value = 55632.668
$ws.Range("L135").Value = 55632.668
$ws.Range("N135").Value = -65772.66800000001
$ws.Range("H136").Value = 27570300
$ws.Range("J136").Value = 1499.6666
$ws.Range("L136").Value = 4498.9998
$ws.Range("N136").Value = -9598.9998
$ws.Range("H138").Value = 49792.31
$ws.Range("J138").Value = 49792.31
$ws.Range("L138").Value = 49792.31
$ws.Range("N138").Value = -60072.31
$ws.Range("H139").Value = 52366.668
$ws.Range("J139").Value = 52366.668
$ws.Range("L139").Value = 52366.668
$ws.Range("N139").Value = -62646.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 77667.62
$ws.Range("I5").Value = 143357.72
$ws.Range("K5").Value = 430073.16
$ws.Range("M5").Value = -429961.16
$ws.Range("H131").Value = 896.36664
$ws.Range("I131").Value = 533.3333
$ws.Range("J131").Value = 987.125
$ws.Range("K131").Value = 1599.9999
$ws.Range("L131").Value = 2961.375
$ws.Range("M131").Value = 3440.0001
$ws.Range("N131").Value = -13041.375
$ws.Range("H132").Value = 2136.762
$ws.Range("I132").Value = 1929
$ws.Range("J132").Value = 2240.6428
$ws.Range("K132").Value = 17361
$ws.Range("L132").Value = 20165.7852
$ws.Range("M132").Value = -14831
$ws.Range("N132").Value = -25225.7852
$ws.Range("H133").Value = 5816.1665
$ws.Range("I133").Value = 4800
$ws.Range("J133").Value = 6832.3335
$ws.Range("K133").Value = 14400
$ws.Range("L133").Value = 20497.0005
$ws.Range("M133").Value = -9340
$ws.Range("N133").Value = -30617.0005
$ws.Range("H135").Value = 77667.62
$ws.Range("I135").Value = 143357.72
$ws.Range("K135").Value = 1290219.48
$ws.Range("M135").Value = -1287684.48

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1236.2727
$ws.Range("I126").Value = 1206.75
$ws.Range("J126").Value = 1315
$ws.Range("K126").Value = 3620.25
$ws.Range("L126").Value = 3945
$ws.Range("M126").Value = -1150.25
$ws.Range("N126").Value = -8885
$ws.Range("H132").Value = 58149.695
$ws.Range("I132").Value = 50097.57
$ws.Range("J132").Value = 69422.664
$ws.Range("K132").Value = 150292.71
$ws.Range("L132").Value = 208267.992
$ws.Range("M132").Value = -147762.71
$ws.Range("N132").Value = -213327.992

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3219.8125
$ws.Range("I61").Value = 2934.55
$ws.Range("J61").Value = 3695.25
$ws.Range("K61").Value = 2934.55
$ws.Range("L61").Value = 3695.25
$ws.Range("M61").Value = -2732.55
$ws.Range("N61").Value = -4099.25
$ws.Range("H113").Value = 3219.8125
$ws.Range("I113").Value = 2934.55
$ws.Range("J113").Value = 3695.25
$ws.Range("K113").Value = 2934.55
$ws.Range("L113").Value = 3695.25
$ws.Range("M113").Value = -764.5500000000002
$ws.Range("N113").Value = -8035.25
$ws.Range("H132").Value = 50386.617
$ws.Range("I132").Value = 1795.6923
$ws.Range("J132").Value = 129346.875
$ws.Range("K132").Value = 5387.0769
$ws.Range("L132").Value = 388040.625
$ws.Range("M132").Value = -2857.0769
$ws.Range("N132").Value = -393100.625
$ws.Range("H133").Value = 29899.54
$ws.Range("J133").Value = 29899.54
$ws.Range("L133").Value = 29899.54
$ws.Range("N133").Value = -34959.54
$ws.Range("H136").Value = 82693.28
$ws.Range("I136").Value = 36170.965
$ws.Range("J136").Value = 275428.56
$ws.Range("K136").Value = 108512.895
$ws.Range("L136").Value = 826285.6799999999
$ws.Range("M136").Value = -105962.895
$ws.Range("N136").Value = -831385.6799999999
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 3360001
$ws.Range("I2").Value = 10000000
$ws.Range("J2").Value = 40001.5
$ws.Range("K2").Value = 10000000
$ws.Range("L2").Value = 40001.5
$ws.Range("M2").Value = -9999888
$ws.Range("N2").Value = -40225.5
$ws.Range("H5").Value = 9999.909
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 9999.909
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 9999.909
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -10223.909
$ws.Range("H123").Value = 41000
$ws.Range("J123").Value = 41000
$ws.Range("L123").Value = 41000
$ws.Range("N123").Value = -50800
$ws.Range("H132").Value = 39521.188
$ws.Range("I132").Value = 27727.105
$ws.Range("K132").Value = 83181.315
$ws.Range("M132").Value = -80651.315
$ws.Range("H135").Value = 48600
$ws.Range("J135").Value = 48600
$ws.Range("L135").Value = 48600
$ws.Range("N135").Value = -58740
$ws.Range("H136").Value = 58964.94
$ws.Range("I136").Value = 34205.6
$ws.Range("J136").Value = 207521
$ws.Range("K136").Value = 102616.8
$ws.Range("L136").Value = 622563
$ws.Range("M136").Value = -100066.8
$ws.Range("N136").Value = -627663
